$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B25").Value = 543.34176935416
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 289296
$ws.Range("E25").Value = 447.2035863348277
$ws.Range("F25").Value = 86788.8
$ws.Range("G25").Value = 72324
$ws.Range("H25").Value = 127
$ws.Range("I25").Value = 98
$ws.Range("J25").Value = 41
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 1233
$ws.Range("O25").Value = 0.29
$ws.Range("P25").Value = 2.4
$ws.Range("Q25").Value = 0.015
$ws.Range("R25").Value = 64000
$ws.Range("S25").Value = 450.8
$ws.Range("T25").Value = 412
$ws.Range("U25").Value = 245
$ws.Range("V25").Value = 170000
$ws.Range("W25").Value = 77000
